$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.021.20"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.008.67"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'225.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("D6").Value = "'0.597"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.84%  "
$ws.Range("D8").Value = "'54.69"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.57%  "
$ws.Range("D9").Value = "'0.374"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.60%  "
$ws.Range("D10").Value = "'0.0780"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("E11").Value = "  -5.35%  "
$ws.Range("D12").Value = "2.306.26"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "'13.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.39%  "
$ws.Range("D14").Value = "'19.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.68%  "
$ws.Range("D15").Value = "'5.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("D16").Value = "'0.734"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.68%  "
$ws.Range("D17").Value = "1.963.10"
$ws.Range("E17").Value = "  -4.33%  "
$ws.Range("D18").Value = "36.945.66"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("E19").Value = "  +2.70%  "
$ws.Range("D20").Value = "'68.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").Value = "0.0₃0811"
$ws.Range("E21").Value = "  -3.67%  "
$ws.Range("D22").Value = "'222.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("D25").Value = "'2.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.21%  "
$ws.Range("D26").Value = "'164.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("D27").Value = "'8.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.61%  "
$ws.Range("E28").Value = "  -4.61%  "
$ws.Range("D29").Value = "'18.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("E30").Value = "  -8.32%  "
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "'4.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.18%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0599"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").Value = "'2.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.28%  "
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("E38").Value = "  -4.60%  "
$ws.Range("D39").Value = "'5.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").Value = "1.452.83"
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("E41").Value = "  -5.01%  "
$ws.Range("D42").Value = "'94.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("D43").Value = "'2.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.48%  "
$ws.Range("D44").Value = "'0.0904"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.64%  "
$ws.Range("D45").Value = "'1.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.51%  "
$ws.Range("D46").Value = "'15.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.87%  "
$ws.Range("D47").Value = "'7.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").Value = "'0.994"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.07%  "
$ws.Range("D49").Value = "'2.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("D50").Value = "2.191.93"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("D51").Value = "'3.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.24%  "
